# daily auto push: 2026-02-23 05:07 UTC
# Insert a new observation row for 2026/02/23 07:00->13:00 slot just above
# the 2026/12/29 block (row 866), shifting the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything currently at/after row 866 down by one row.
$ws.Rows(866).Insert()

# Fill in the newly inserted row with the new day's data.
# Column A holds a date-like string ("2026/02/23") that must stay literal
# text (matching the rest of the sheet) instead of being auto-converted
# to a real Excel date serial number. Temporarily mark the cell as Text,
# assign the value, then clear the formatting back to General so the
# saved cell carries no special number format / style, just like its
# neighbours.
$ws.Range("A866").NumberFormat = "@"
$ws.Range("A866").Value = "2026/02/23"
$ws.Range("A866").ClearFormats()

$ws.Range("B866").Value = "月"
$ws.Range("C866").Value = 13
$ws.Range("D866").Value = 201
